# Scheduled market-data refresh: update computed Leve profit columns (H:N)
# on each job sheet with the latest Universalis price snapshot values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2: Mercury Rising
$ws.Range("H2").Value = 65
$ws.Range("J2").Value = 100
$ws.Range("L2").Value = 100
$ws.Range("N2").Value = -326
# Row 38: Just Give Him a Serum
$ws.Range("H38").Value = 558.8
$ws.Range("I38").Value = 136
$ws.Range("J38").Value = 2250
$ws.Range("K38").Value = 408
$ws.Range("L38").Value = 6750
$ws.Range("M38").Value = -36
$ws.Range("N38").Value = -7494
# Row 58: A Matter of Vital Importance
$ws.Range("H58").Value = 2412.8572
$ws.Range("J58").Value = 4714.2856
$ws.Range("L58").Value = 14142.8568
$ws.Range("N58").Value = -14442.8568
# Row 98: The Dotted Line
$ws.Range("H98").Value = 776.15
$ws.Range("I98").Value = 734.6667
$ws.Range("J98").Value = 1149.5
$ws.Range("K98").Value = 734.6667
$ws.Range("L98").Value = 1149.5
$ws.Range("M98").Value = 763.3333
$ws.Range("N98").Value = -4145.5
# Row 122: Wishful Inking
$ws.Range("H122").Value = 776.15
$ws.Range("I122").Value = 734.6667
$ws.Range("J122").Value = 1149.5
$ws.Range("K122").Value = 2204.0001
$ws.Range("L122").Value = 3448.5
$ws.Range("M122").Value = 245.9998999999998
$ws.Range("N122").Value = -8348.5
# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 4549175
$ws.Range("I137").Value = 5558159
$ws.Range("J137").Value = 8747.5
$ws.Range("K137").Value = 16674477
$ws.Range("L137").Value = 26242.5
$ws.Range("M137").Value = -16671927
$ws.Range("N137").Value = -31342.5
# Row 138: All-night Crafting
$ws.Range("H138").Value = 2010432.6
$ws.Range("I138").Value = 690.5
$ws.Range("J138").Value = 3271447.2
$ws.Range("K138").Value = 2071.5
$ws.Range("L138").Value = 9814341.600000001
$ws.Range("M138").Value = 3068.5
$ws.Range("N138").Value = -9824621.600000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 14473.429
$ws.Range("I32").Value = 16342.118
$ws.Range("J32").Value = 8017.9546
$ws.Range("K32").Value = 16342.118
$ws.Range("L32").Value = 8017.9546
$ws.Range("M32").Value = -16055.118
$ws.Range("N32").Value = -8591.954600000001
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 143144290
$ws.Range("I61").Value = 250250750
$ws.Range("J61").Value = 335671.34
$ws.Range("K61").Value = 250250750
$ws.Range("L61").Value = 335671.34
$ws.Range("M61").Value = -250250538
$ws.Range("N61").Value = -336095.34
# Row 122: Haste for High Durium
$ws.Range("H122").Value = 3004796
$ws.Range("I122").Value = 1833.5
$ws.Range("J122").Value = 37038372
$ws.Range("K122").Value = 5500.5
$ws.Range("L122").Value = 111115116
$ws.Range("M122").Value = -3050.5
$ws.Range("N122").Value = -111120016
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 143144290
$ws.Range("I136").Value = 250250750
$ws.Range("J136").Value = 335671.34
$ws.Range("K136").Value = 750752250
$ws.Range("L136").Value = 1007014.02
$ws.Range("M136").Value = -750749700
$ws.Range("N136").Value = -1012114.02

$ws = $wb.Worksheets.Item("BSM")
# Row 19: Twice as Slice
$ws.Range("H19").Value = 6100
$ws.Range("J19").Value = 6100
$ws.Range("L19").Value = 6100
$ws.Range("N19").Value = -6446
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 7839.6
$ws.Range("I134").Value = 6924.5
$ws.Range("J134").Value = 11500
$ws.Range("K134").Value = 20773.5
$ws.Range("L134").Value = 34500
$ws.Range("M134").Value = -18238.5
$ws.Range("N134").Value = -39570

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 125000664
$ws.Range("I22").Value = 142857800
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 142857800
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = -142857450
$ws.Range("N22").Value = -1400
# Row 31: Wall Not Found
$ws.Range("H31").Value = 3426.375
$ws.Range("I31").Value = 2218.7222
$ws.Range("J31").Value = 4979.0713
$ws.Range("K31").Value = 2218.7222
$ws.Range("L31").Value = 4979.0713
$ws.Range("M31").Value = -1923.7222
$ws.Range("N31").Value = -5569.0713
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 3426.375
$ws.Range("I34").Value = 2218.7222
$ws.Range("J34").Value = 4979.0713
$ws.Range("K34").Value = 2218.7222
$ws.Range("L34").Value = 4979.0713
$ws.Range("M34").Value = -2016.7222
$ws.Range("N34").Value = -5383.0713
# Row 99: O Pine
$ws.Range("H99").Value = 3363.889
$ws.Range("I99").Value = 2601.3157
$ws.Range("J99").Value = 5175
$ws.Range("K99").Value = 2601.3157
$ws.Range("L99").Value = 5175
$ws.Range("M99").Value = -1103.3157
$ws.Range("N99").Value = -8171
# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 927.5
$ws.Range("I105").Value = 872.8333
$ws.Range("J105").Value = 1255.5
$ws.Range("K105").Value = 872.8333
$ws.Range("L105").Value = 1255.5
$ws.Range("M105").Value = 874.1667
$ws.Range("N105").Value = -4749.5
# Row 126: A Better Conductor
$ws.Range("H126").Value = 3363.889
$ws.Range("I126").Value = 2601.3157
$ws.Range("J126").Value = 5175
$ws.Range("K126").Value = 7803.9471
$ws.Range("L126").Value = 15525
$ws.Range("M126").Value = -5333.9471
$ws.Range("N126").Value = -20465

$ws = $wb.Worksheets.Item("CUL")
# Row 86: Let's Not Get Sappy
$ws.Range("H86").Value = 1020.6
$ws.Range("J86").Value = 1367.6666
$ws.Range("L86").Value = 4102.9998
$ws.Range("N86").Value = -6474.9998
# Row 89: Luxury Spillover (L)
$ws.Range("H89").Value = 1020.6
$ws.Range("J89").Value = 1367.6666
$ws.Range("L89").Value = 12308.9994
$ws.Range("N89").Value = -24164.9994

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 57955.26
$ws.Range("I70").Value = 115861.11
$ws.Range("J70").Value = 5840
$ws.Range("K70").Value = 115861.11
$ws.Range("L70").Value = 5840
$ws.Range("M70").Value = -115591.11
$ws.Range("N70").Value = -6380
# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 57955.26
$ws.Range("I73").Value = 115861.11
$ws.Range("J73").Value = 5840
$ws.Range("K73").Value = 115861.11
$ws.Range("L73").Value = 5840
$ws.Range("M73").Value = -114925.11
$ws.Range("N73").Value = -7712
# Row 123: Workplace Workout
$ws.Range("H123").Value = 38468
$ws.Range("J123").Value = 38468
$ws.Range("L123").Value = 38468
$ws.Range("N123").Value = -43368

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 41747.26
$ws.Range("I132").Value = 19048.035
$ws.Range("J132").Value = 114889.22
$ws.Range("K132").Value = 57144.105
$ws.Range("L132").Value = 344667.66
$ws.Range("M132").Value = -54614.105
$ws.Range("N132").Value = -349727.66
# Row 133: The Perfect Accessory
$ws.Range("H133").Value = 37329.223
$ws.Range("J133").Value = 38870.375
$ws.Range("L133").Value = 38870.375
$ws.Range("N133").Value = -43930.375
# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 46786.312
$ws.Range("I136").Value = 29008.445
$ws.Range("J136").Value = 117897.78
$ws.Range("K136").Value = 87025.33499999999
$ws.Range("L136").Value = 353693.34
$ws.Range("M136").Value = -84475.33499999999
$ws.Range("N136").Value = -358793.34

$ws = $wb.Worksheets.Item("WVR")
# Row 126: A Polished Purchase
$ws.Range("H126").Value = 1412.4667
$ws.Range("I126").Value = 1258.7
$ws.Range("J126").Value = 1720
$ws.Range("K126").Value = 3776.1
$ws.Range("L126").Value = 5160
$ws.Range("M126").Value = -1306.1
$ws.Range("N126").Value = -10100
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 78346.16
$ws.Range("I132").Value = 50950.1
$ws.Range("J132").Value = 169666.33
$ws.Range("K132").Value = 152850.3
$ws.Range("L132").Value = 508998.99
$ws.Range("M132").Value = -150320.3
$ws.Range("N132").Value = -514058.99
# Row 135: In Line with Linen
$ws.Range("H135").Value = 33700
$ws.Range("J135").Value = 33700
$ws.Range("L135").Value = 33700
$ws.Range("N135").Value = -43840
# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 40203.5
$ws.Range("I136").Value = 26968.77
$ws.Range("K136").Value = 80906.31
$ws.Range("M136").Value = -78356.31
